# Update rules in DiscountRules.xlsx
#
# The "RuleTable" block on the sheet (header in row 18, rule rows 19-26)
# is widened from 3 columns (NAME | CONDITION | ACTION) to 5 columns
# (NAME | CONDITION | CONDITION | ACTION | ACTION): a second CONDITION
# column is inserted after the existing one, and a second ACTION column
# is appended after the existing ACTION column. The data that used to live
# in column C (the single ACTION/condition value per row) slides over into
# the new column D, while the freshly inserted columns C & E are populated
# with the new rule note, or left blank for rows that don't carry a note.
#
# A new code-change note "Code changed 10010018" supersedes the previous
# "Code changed 10010017" note, and a new trailing row (26) is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newNote = "Code changed 10010018"

function Set-BlankCell($sheet, $row, $col) {
    # Produces a present-but-empty cell (rather than leaving the cell out
    # of sheetData entirely) by touching its style after clearing the
    # value - mirrors the blank placeholder cells the rule-table exporter
    # leaves behind when a column has no text for a given row.
    $sheet.Cells.Item($row, $col).Value = ""
    $sheet.Cells.Item($row, $col).Style = "Normal"
}

# --- Row 18 (header row: NAME | CONDITION | ACTION) -----------------------
# Duplicate the CONDITION header into the new column C and push the
# existing ACTION header out to columns D & E.
$ws.Cells.Item(18, 4).Value = $ws.Cells.Item(18, 3).Value()
$ws.Cells.Item(18, 3).Value = $ws.Cells.Item(18, 2).Value()
$ws.Cells.Item(18, 5).Value = $ws.Cells.Item(18, 4).Value()

# --- Row 19 (CONDITION expression | ACTION "Test") -------------------------
$ws.Cells.Item(19, 4).Value = $ws.Cells.Item(19, 3).Value()
$ws.Cells.Item(19, 3).Value = $newNote
$ws.Cells.Item(19, 5).Value = $newNote

# --- Rows 20-23 (rule rows: shift old column C into D, blank out C & E) ---
for ($r = 20; $r -le 23; $r++) {
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 3).Value()
    Set-BlankCell $ws $r 3
    Set-BlankCell $ws $r 5
}

# --- Row 24 (blank ACTION cell already; keep C/D/E blank) -----------------
Set-BlankCell $ws 24 3
Set-BlankCell $ws 24 4
Set-BlankCell $ws 24 5

# --- Row 25 (code-change note row, bumped to the new note) -----------------
$ws.Cells.Item(25, 1).Value = $newNote
$ws.Cells.Item(25, 3).Value = $newNote
Set-BlankCell $ws 25 4
Set-BlankCell $ws 25 5

# --- Row 26 (new trailing row) ----------------------------------------------
Set-BlankCell $ws 26 1
Set-BlankCell $ws 26 2
$ws.Cells.Item(26, 3).Value = $newNote
Set-BlankCell $ws 26 4
$ws.Cells.Item(26, 5).Value = $newNote
